$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 currently has the "bottom of table" border style; once its content
# comes from row 21 (which uses the regular interior style), copy row 21's
# formatting onto D22 so it matches the rest of the shifted column.
$ws.Range("D21").Copy()
$ws.Range("D22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Shift the "Ferreteria" (column D) values down by one row to make room for
# the new "Terminales Dupont" entry, then set the new row's value.
$ws.Range("D22").Value = $ws.Range("D21").Value2
$ws.Range("D21").Value = $ws.Range("D20").Value2
$ws.Range("D20").Value = $ws.Range("D19").Value2
$ws.Range("D19").Value = $ws.Range("D18").Value2
$ws.Range("D18").Value = "Terminales Dupont"

$ws.Range("E18").Select()
